$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 85, shifting existing rows 85:90 down to 86:91.
$ws.Rows.Item(85).Insert()

# Populate the newly inserted row 85 with a new weekly price record, matching the
# same constant columns used by the surrounding "Vega Modelo de Temuco" / "Arveja Verde" rows.
$ws.Cells.Item(85, 1).Value = 10
$ws.Cells.Item(85, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(85, 3).Value = "La Araucanía"
$ws.Cells.Item(85, 4).Value = 44585
$ws.Cells.Item(85, 5).Value = 9
$ws.Cells.Item(85, 6).Value = 100112022
$ws.Cells.Item(85, 7).Value = "Arveja Verde"
$ws.Cells.Item(85, 8).Value = "Sin especificar"
$ws.Cells.Item(85, 9).Value = "Primera"
$ws.Cells.Item(85, 10).Value = 65
$ws.Cells.Item(85, 11).Value = 26000
$ws.Cells.Item(85, 12).Value = 26000
$ws.Cells.Item(85, 13).Value = 26000
$ws.Cells.Item(85, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(85, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(85, 16).Value = 1040
$ws.Cells.Item(85, 17).Value = 25
$ws.Cells.Item(85, 18).Value = "Hortaliza"
